$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sqlCount")

# The sqlRecordCount/sqlColCount result row (row 2) needs to be updated to
# reflect the corrected query results: 209 records, 0 columns.
# Force the cells to stay text (matching the original "252"/"1" text cells)
# instead of being auto-converted to numbers, then restore General format.
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "209"
$ws.Range("B2").Value = "0"
$ws.Range("A2:B2").NumberFormat = "General"
